$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Replacement.ClearFormatting()
    $ok = $range.Find.Execute($find, $true, $true, $false, $false, $false, `
                               $true, 1, $false, $replace, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $find"
    }
}

# Change 1: "Mayornan/Edukadónan por inisialmente ke" -> "Mayornan/Edukadónan inisialmente lo ke"
Replace-Text "Mayornan/Edukadónan por inisialmente ke" "Mayornan/Edukadónan inisialmente lo ke"

# Change 2: "un mucha, i tene e mesun tono" -> "un mucha, Tene e mesun tono"
Replace-Text "un mucha, i tene e mesun tono" "un mucha, Tene e mesun tono"

# Change 3: "no ta kustumá ku ta puntra" -> "no ta kustumbrá ku ta puntra"
Replace-Text "no ta kustumá ku ta puntra" "no ta kustumbrá ku ta puntra"

# Change 3b: "no ta kustuma ku hende" -> "no ta kustumbrá ku hende"
Replace-Text "no ta kustuma ku hende" "no ta kustumbrá ku hende"

# Change 4: "e kosnan ei ku nan no por hasi" -> "e kosnan ku nan no por hasi"
Replace-Text "e kosnan ei ku nan no por hasi" "e kosnan ku nan no por hasi"

# Change 4b: "bon den pintamentu. " -> "bon den pintamentu"
Replace-Text "bon den pintamentu. " "bon den pintamentu"
